# "Generate Report for Handback" — update the localization-status report
# after a handback transform failure: the status moves from
# "Ready for handoff" to "Handback transform failed", and the per-language
# Error Detail column is populated with the failure detail. Columns that
# now hold longer text are widened to fit.

$wb = $excel.ActiveWorkbook

$statusNew = "Handback transform failed"
$errorDetail = "The translationStateItem 6cfb131bc7167d1c9d969dbc3dfc212f8169e884 is not found."

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: zh-cn / de-de status columns (E2, F2)
$overview.Range("E2").Value = $statusNew
$overview.Range("F2").Value = $statusNew

# zh-cn sheet: Status (C2) + Error Detail (P2)
$zhcn.Range("C2").Value = $statusNew
$zhcn.Range("P2").Value = $errorDetail

# de-de sheet: Status (C2) + Error Detail (P2)
$dede.Range("C2").Value = $statusNew
$dede.Range("P2").Value = $errorDetail

# Widen columns to fit the new, longer text.
$overview.Columns.Item(5).ColumnWidth = 23.75
$overview.Columns.Item(6).ColumnWidth = 23.75

$zhcn.Columns.Item(3).ColumnWidth = 23.75
$zhcn.Columns.Item(16).ColumnWidth = 39.16

$dede.Columns.Item(3).ColumnWidth = 23.75
$dede.Columns.Item(16).ColumnWidth = 39.16
